# Applies the cryptocurrency price/volume/ranking update described in the commit
# message "Updated cryptos list on Thu Feb 29 20:53:03 UTC 2024 with GitHub Actions".
#
# All data cells in this sheet are stored as TEXT (inline strings) in the source
# workbook, including the "Price" column which often contains plain-looking numeric
# literals (e.g. "407.02"). Assigning such a literal straight to Range.Value makes
# Excel auto-detect it as a number and restyle the cell, which would not match the
# original text-only layout. To avoid that, for any new value that looks like a bare
# number we first force the cells NumberFormat to Text ("@") before assigning the
# value, then call ClearFormats() to drop the incidental style Excel attached, so the
# cell ends up back in its default (style-less) state but still holding a text value.

function Set-CellText($range, [string]$text) {
    if ($text -match '^-?\d+(\.\d+)?$') {
        $range.NumberFormat = "@"
        $range.Value = $text
        $range.ClearFormats()
    } else {
        $range.Value = $text
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-CellText $ws.Range("D2") "62.176.79"
Set-CellText $ws.Range("E2") "  +3.20%  "
# Row 3
Set-CellText $ws.Range("D3") "3.406.91"
Set-CellText $ws.Range("E3") "  +3.46%  "
# Row 4
Set-CellText $ws.Range("E4") "  +0.10%  "
# Row 5
Set-CellText $ws.Range("D5") "407.02"
Set-CellText $ws.Range("E5") "  -0.29%  "
# Row 6
Set-CellText $ws.Range("D6") "130.84"
Set-CellText $ws.Range("E6") "  +16.70%  "
# Row 7
Set-CellText $ws.Range("E7") "  +6.90%  "
# Row 8
Set-CellText $ws.Range("E8") "  -0.04%  "
# Row 9
Set-CellText $ws.Range("D9") "0.678"
Set-CellText $ws.Range("E9") "  +9.30%  "
# Row 10
Set-CellText $ws.Range("D10") "0.127"
Set-CellText $ws.Range("E10") "  +10.85%  "
# Row 11
Set-CellText $ws.Range("D11") "42.22"
Set-CellText $ws.Range("E11") "  +8.85%  "
# Row 12
Set-CellText $ws.Range("E12") "  -0.39%  "
# Row 13
Set-CellText $ws.Range("D13") "3.962.23"
Set-CellText $ws.Range("E13") "  +4.21%  "
# Row 14
Set-CellText $ws.Range("D14") "8.58"
Set-CellText $ws.Range("E14") "  +5.06%  "
# Row 15
Set-CellText $ws.Range("E15") "  +4.28%  "
# Row 16
Set-CellText $ws.Range("D16") "3.403.08"
Set-CellText $ws.Range("E16") "  +2.53%  "
# Row 17
Set-CellText $ws.Range("B17") "Uniswap"
Set-CellText $ws.Range("C17") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-CellText $ws.Range("D17") "11.58"
Set-CellText $ws.Range("E17") "  +8.15%  "
# Row 18
Set-CellText $ws.Range("B18") "WrappedBTC"
Set-CellText $ws.Range("C18") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-CellText $ws.Range("D18") "62.023.62"
Set-CellText $ws.Range("E18") "  +3.35%  "
# Row 19
Set-CellText $ws.Range("E19") "  +4.85%  "
# Row 20
Set-CellText $ws.Range("D20") "0.0000134"
Set-CellText $ws.Range("E20") "  +17.57%  "
# Row 21
Set-CellText $ws.Range("E21") "  -0.59%  "
# Row 22
Set-CellText $ws.Range("D22") "82.76"
Set-CellText $ws.Range("E22") "  +12.57%  "
# Row 23
Set-CellText $ws.Range("D23") "13.21"
Set-CellText $ws.Range("E23") "  +6.14%  "
# Row 24
Set-CellText $ws.Range("D24") "308.26"
Set-CellText $ws.Range("E24") "  +4.25%  "
# Row 25
Set-CellText $ws.Range("E25") "  +2.52%  "
# Row 26
Set-CellText $ws.Range("D26") "8.61"
Set-CellText $ws.Range("E26") "  +14.82%  "
# Row 27
Set-CellText $ws.Range("B27") "EthereumClassic"
Set-CellText $ws.Range("C27") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-CellText $ws.Range("D27") "29.79"
Set-CellText $ws.Range("E27") "  +2.24%  "
# Row 28
Set-CellText $ws.Range("B28") "LEO"
Set-CellText $ws.Range("C28") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-CellText $ws.Range("D28") "4.68"
Set-CellText $ws.Range("E28") "  +9.71%  "
# Row 29
Set-CellText $ws.Range("D29") "0.175"
Set-CellText $ws.Range("E29") "  +1.90%  "
# Row 30
Set-CellText $ws.Range("D30") "7.48"
Set-CellText $ws.Range("E30") "  +1.82%  "
# Row 31
Set-CellText $ws.Range("E31") "  +2.58%  "
# Row 32
Set-CellText $ws.Range("D32") "11.80"
Set-CellText $ws.Range("E32") "  +5.53%  "
# Row 33
Set-CellText $ws.Range("E33") "  +6.25%  "
# Row 34
Set-CellText $ws.Range("D34") "42.49"
Set-CellText $ws.Range("E34") "  +8.56%  "
# Row 35
Set-CellText $ws.Range("D35") "1.00"
Set-CellText $ws.Range("E35") "  -0.05%  "
# Row 36
Set-CellText $ws.Range("E36") "  +1.58%  "
# Row 37
Set-CellText $ws.Range("D37") "52.45"
Set-CellText $ws.Range("E37") "  +0.88%  "
# Row 38
Set-CellText $ws.Range("D38") "0.999"
Set-CellText $ws.Range("E38") "  +0.25%  "
# Row 39
Set-CellText $ws.Range("E39") "  +3.97%  "
# Row 40
Set-CellText $ws.Range("E40") "  -2.92%  "
# Row 41
Set-CellText $ws.Range("E41") "  +8.77%  "
# Row 42
Set-CellText $ws.Range("E42") "  +5.15%  "
# Row 43
Set-CellText $ws.Range("D43") "137.63"
Set-CellText $ws.Range("E43") "  +2.52%  "
# Row 44
Set-CellText $ws.Range("D44") "3.99"
Set-CellText $ws.Range("E44") "  +5.49%  "
# Row 45
Set-CellText $ws.Range("B45") "Celestia"
Set-CellText $ws.Range("C45") "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-CellText $ws.Range("D45") "17.10"
Set-CellText $ws.Range("E45") "  +5.30%  "
# Row 46
Set-CellText $ws.Range("B46") "TheGraph"
Set-CellText $ws.Range("C46") "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-CellText $ws.Range("D46") "0.286"
Set-CellText $ws.Range("E46") "  -2.92%  "
# Row 47
Set-CellText $ws.Range("E47") "  +2.15%  "
# Row 48
Set-CellText $ws.Range("D48") "21.74"
Set-CellText $ws.Range("E48") "  +4.49%  "
# Row 49
Set-CellText $ws.Range("B49") "Maker"
Set-CellText $ws.Range("C49") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-CellText $ws.Range("D49") "2.152.20"
Set-CellText $ws.Range("E49") "  +1.38%  "
# Row 50
Set-CellText $ws.Range("B50") "RocketPoolETH"
Set-CellText $ws.Range("C50") "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-CellText $ws.Range("D50") "3.750.20"
Set-CellText $ws.Range("E50") "  +4.03%  "
# Row 51
Set-CellText $ws.Range("E51") "  -0.76%  "
